$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAEPAbCiPC")
$ws.Range("B22").Value = 0
$ws.Range("D26").Select() | Out-Null
$wb.Worksheets.Item("About").Activate() | Out-Null
